# Semana 24 de 2025: add columns Z (week 23) and AA (week 24)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new week-number labels, bold + centered like the rest of the header ---
$ws.Range("Z1").Value = "'23"
$ws.Range("Z1").Font.Bold = $true
$ws.Range("Z1").HorizontalAlignment = -4108

$ws.Range("AA1").Value = "'24"
$ws.Range("AA1").Font.Bold = $true
$ws.Range("AA1").HorizontalAlignment = -4108

# --- Data rows: weekly case counts for week 23 (Z) and week 24 (AA) ---
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = 0
$ws.Range("Z3").Value = 0
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 0
$ws.Range("Z5").Value = 0
$ws.Range("AA5").Value = 0
$ws.Range("Z6").Value = 2
$ws.Range("AA6").Value = 1
$ws.Range("Z7").Value = 0
$ws.Range("AA7").Value = 0
$ws.Range("Z8").Value = 0
$ws.Range("AA8").Value = 0
$ws.Range("Z9").Value = 0
$ws.Range("AA9").Value = 0
$ws.Range("Z10").Value = 0
$ws.Range("Z11").Value = 0
$ws.Range("Z12").Value = 0
$ws.Range("AA12").Value = 0
$ws.Range("Z13").Value = 0
$ws.Range("AA13").Value = 0
$ws.Range("Z14").Value = 0
$ws.Range("AA14").Value = 0
$ws.Range("Z15").Value = 0
$ws.Range("AA15").Value = 0
$ws.Range("Z16").Value = 0
$ws.Range("AA16").Value = 0
$ws.Range("Z17").Value = 0
$ws.Range("AA17").Value = 0
$ws.Range("Z18").Value = 0
$ws.Range("Z19").Value = 0
$ws.Range("Z21").Value = 0
$ws.Range("AA21").Value = 0
$ws.Range("Z22").Value = 0
$ws.Range("AA22").Value = 0
$ws.Range("AA23").Value = 0
$ws.Range("Z24").Value = 0
$ws.Range("AA24").Value = 0
$ws.Range("Z25").Value = 0
$ws.Range("AA25").Value = 0
$ws.Range("Z27").Value = 15
$ws.Range("AA27").Value = 0
$ws.Range("Y28").Value = 1
$ws.Range("Z28").Value = 1
$ws.Range("AA28").Value = 1
$ws.Range("Y29").Value = 2
$ws.Range("Z29").Value = 5
$ws.Range("AA29").Value = 3
$ws.Range("Z30").Value = 0
$ws.Range("AA30").Value = 0
$ws.Range("Z31").Value = 0
$ws.Range("AA31").Value = 0
$ws.Range("Z33").Value = 0
$ws.Range("AA33").Value = 0
$ws.Range("Z34").Value = 3
$ws.Range("AA34").Value = 5
$ws.Range("Y35").Value = 0
$ws.Range("Z35").Value = 0
$ws.Range("AA35").Value = 0
$ws.Range("Z36").Value = 0
$ws.Range("AA36").Value = 0
$ws.Range("Z37").Value = 0
$ws.Range("AA37").Value = 0
$ws.Range("Z39").Value = 0
$ws.Range("AA39").Value = 0
$ws.Range("Z40").Value = 0
$ws.Range("Z41").Value = 0
$ws.Range("Z42").Value = 0
$ws.Range("AA42").Value = 0
$ws.Range("Z43").Value = 0
$ws.Range("AA43").Value = 0
$ws.Range("Z44").Value = 0
$ws.Range("AA44").Value = 0
$ws.Range("Z45").Value = 0
$ws.Range("AA45").Value = 0
$ws.Range("Z46").Value = 0
$ws.Range("AA46").Value = 0
$ws.Range("Z47").Value = 0
$ws.Range("AA47").Value = 0
$ws.Range("Z48").Value = 0
$ws.Range("AA48").Value = 0
$ws.Range("AA49").Value = 0
$ws.Range("Y50").Value = 0
$ws.Range("Z50").Value = 0
$ws.Range("Z51").Value = 0
$ws.Range("AA51").Value = 0
$ws.Range("Z52").Value = 0
$ws.Range("AA52").Value = 0
$ws.Range("Z53").Value = 0
$ws.Range("AA53").Value = 0
$ws.Range("Z54").Value = 0
$ws.Range("AA54").Value = 0
$ws.Range("Z55").Value = 0
$ws.Range("AA55").Value = 0
$ws.Range("Z56").Value = 0
$ws.Range("AA56").Value = 0
